$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column (C) for rows 2-8
# from serial date 45207 (2023-10-08) to 45208 (2023-10-09)
$ws.Range("C2:C8").Value = 45208
